$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect so values can be written, then
# restore protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A41).
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns for every holding row.
$ws.Range("D2").Value = 0.03205282870608252
$ws.Range("E2").Value = 0.003950871768444664
$ws.Range("D3").Value = 0.02870325088385083
$ws.Range("E3").Value = -0.007851035843472576
$ws.Range("D4").Value = 0.02786811779643637
$ws.Range("E4").Value = 0.0008043973723019793
$ws.Range("D5").Value = 0.06337533382218215
$ws.Range("E5").Value = 0.01473517639393096
$ws.Range("D6").Value = 0.01583115012987136
$ws.Range("E6").Value = -0.004098971530779649
$ws.Range("D7").Value = 0.01535567289762079
$ws.Range("E7").Value = -0.002689204763734132
$ws.Range("D8").Value = 0.02965892472534701
$ws.Range("E8").Value = 0.003978041212507
$ws.Range("D9").Value = 0.03440425829386264
$ws.Range("E9").Value = 0.003429355281207025
$ws.Range("D10").Value = 0.02910498785055386
$ws.Range("E10").Value = -0.008979062366986312
$ws.Range("D11").Value = 0.02883873633175018
$ws.Range("E11").Value = 0.007350484463748863
$ws.Range("D12").Value = 0.01114166252246374
$ws.Range("E12").Value = -0.01526650194140466
$ws.Range("D13").Value = 0.01424386633387193
$ws.Range("E13").Value = 0.006184770003865525
$ws.Range("D14").Value = 0.0144131739835691
$ws.Range("E14").Value = -0.02292044694871531
$ws.Range("D15").Value = 0.009188627009348497
$ws.Range("E15").Value = 0.01121383324773162
$ws.Range("D16").Value = 0.008282506626300264
$ws.Range("E16").Value = 0.002350427350427342
$ws.Range("D17").Value = 0.0299605715716135
$ws.Range("E17").Value = -0.01194523568868866
$ws.Range("D18").Value = 0.02546123552896449
$ws.Range("E18").Value = -0.007136182142553515
$ws.Range("D19").Value = 0.0333777937974408
$ws.Range("E19").Value = 0.002916224814422108
$ws.Range("D20").Value = 0.03106333266278852
$ws.Range("E20").Value = -0.001519275811863063
$ws.Range("D21").Value = 0.04554513425349387
$ws.Range("E21").Value = 0.002266682785299734
$ws.Range("D22").Value = 0.0362575969679576
$ws.Range("E22").Value = 0.001301624318680972
$ws.Range("D23").Value = 0.03181941618553602
$ws.Range("E23").Value = -0.01118561320025968
$ws.Range("D24").Value = 0.03123460671956696
$ws.Range("E24").Value = -0.002379738228795003
$ws.Range("D25").Value = 0.0149138212264599
$ws.Range("E25").Value = 0.003203987184051238
$ws.Range("D26").Value = 0.01494508709905206
$ws.Range("E26").Value = -0.01531538643719899
$ws.Range("D27").Value = 0.03129497541696188
$ws.Range("E27").Value = 0.002972076306331228
$ws.Range("D28").Value = 0.03067595046777875
$ws.Range("E28").Value = 0.02235897435897449
$ws.Range("D29").Value = 0.0292778350329973
$ws.Range("E29").Value = -0.01196856739875085
$ws.Range("D30").Value = 0.02940230860124156
$ws.Range("E30").Value = 0.003531229309203354
$ws.Range("D31").Value = 0.03360904327019791
$ws.Range("E31").Value = -0.005441269395492498
$ws.Range("D32").Value = 0.03170811754473626
$ws.Range("E32").Value = -0.00286513404734301
$ws.Range("D33").Value = 0.02874847824671999
$ws.Range("E33").Value = -0.02365285434821285
$ws.Range("D34").Value = 0.03223334487614291
$ws.Range("E34").Value = -0.002855051244509577
$ws.Range("D35").Value = 0.03053908853492251
$ws.Range("D36").Value = 0.03225104253987431
$ws.Range("E36").Value = 0.004024144869215318
$ws.Range("D37").Value = 0.03321812154244185
$ws.Range("E37").Value = -0.00539874976321264
$ws.Range("E38").Value = -0.0008652191157577649

# Restore sheet protection.
$ws.Protect()
